# Sprint 3 Backlog - Burndown.xlsx update
# - Add two new "filter recipes" UI tasks (rows 9 & 10), shifting the
#   existing "add to planned meals" / "view different pages of recipes" /
#   "Fix testing problems" rows down to rows 11-15.
# - Fill in Actual Time / Week 1 / Week 2 progress numbers for several tasks.
# - Insert two extra blank rows above the totals block (old row 22 becomes
#   row 24) and refresh the SUM/SUMIF formulas to cover the new data range.
# - Point the burndown chart + the hidden _FilterDatabase name at the
#   relocated totals row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update/extend the existing task rows (still at their original row
#     numbers; the totals block below is inserted afterwards). ---------

# Row 4 - Create UI for planned meals page(web): log actual time worked
$ws.Range("D4").Value = 1
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0

# Row 6 - Implement functionality for view planned meals(Web)
$ws.Range("C6").Value = 2.5
$ws.Range("D6").Value = 6
$ws.Range("F6").Value = 2.5
$ws.Range("G6").Value = 0

# Row 7 - Implement functionality to filter recipes by tags(Web)
$ws.Range("D7").Value = 5
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0

# Row 9 - replaced with a new filter-recipes UI task
$ws.Range("A9").Value = "I want to be able to filter recipes"
$ws.Range("B9").Value = "Modify Recipe page UI to use checkbox and button instead of radio buttons"
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = "Destiny"
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0

# Row 10 - replaced with a new filter-recipes functionality task
$ws.Range("A10").Value = "I want to be able to filter recipes"
$ws.Range("B10").Value = "Modify functionality to fit new UI for filtering recipes on Recipe page"
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 4
$ws.Range("E10").Value = "Destiny"
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = 0

# Row 11 - (previously row 9) add-to-planned-meals desktop task
$ws.Range("A11").Value = "I want to be able to add to planned meals"
$ws.Range("B11").Value = "Implement functionality to add recipe to a specific day of the week for planned meals(desktop)"
$ws.Range("C11").Value = 1.5
$ws.Range("E11").Value = "Janera"

# Row 12 - (previously row 10) add-to-planned-meals web task
$ws.Range("A12").Value = "I want to be able to add to planned meals"
$ws.Range("B12").Value = "Implement functionality to add recipe to a specific day of the week for planned meals(web)"
$ws.Range("C12").Value = 1.5
$ws.Range("D12").Value = 1.5
$ws.Range("E12").Value = "Matthew"
$ws.Range("F12").Value = 2
$ws.Range("G12").Value = 0

# Row 13 - (previously row 11) view-different-pages desktop task
$ws.Range("A13").Value = "I want to be able to view different pages of recipes"
$ws.Range("B13").Value = "Implement functionality to make recipe list paginated(Desktop)"
$ws.Range("C13").Value = 1.5
$ws.Range("E13").Value = "Janera"

# Row 14 - (previously row 12) view-different-pages web task
$ws.Range("A14").Value = "I want to be able to view different pages of recipes"
$ws.Range("B14").Value = "Implement functionality to make recipe list paginated(Web)"
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 6
$ws.Range("E14").Value = "Matthew"
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0

# Row 15 - (previously row 13) Fix testing problems
$ws.Range("B15").Value = "Fix testing problems"
$ws.Range("C15").Value = 3
$ws.Range("E15").Value = "Destiny"
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 3

# --- Make room for two more blank rows above the totals block ----------
$ws.Rows("22:23").Insert()

# --- Refresh the totals formulas to span the new data range (rows 3-23) -
$ws.Range("C24").Formula = "=SUM(C5:C23)"
$ws.Range("D24").Formula = "=SUM(D3:D23)"
$ws.Range("F24").Formula = "=SUM(F3:F23)"
$ws.Range("G24").Formula = "=SUM(G3:G23)"

$ws.Range("G25").Formula = '=SUMIF(E3:E23, "Matthew", C3:C23)'
$ws.Range("G26").Formula = '=SUMIF(E3:E23, "Destiny", C3:C23)'
$ws.Range("G27").Formula = '=SUMIF(E3:E23,"Janera",C3:C23)'

# --- Point the hidden AutoFilter-database name at the relocated totals --
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Sheet1!`$B`$24:`$G`$27"
    }
}

# --- Re-point the burndown chart series at the relocated totals row ----
$chart = $ws.ChartObjects(1)
$co = $chart.Chart
$series = $co.SeriesCollection(1)
$series.Formula = "=SERIES(,,Sheet1!`$F`$24:`$G`$24,1)"

# Grow the chart so its anchor keeps pace with the two inserted rows
$chart.Height = $chart.Height + 30

# --- Restore the active selection ---------------------------------------
$ws.Range("D21").Select()
